# 197-MS-EI-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME-Loanproduct.xlsx
#
# The product "short code" string that lives in cell B1 of both
# ProductLoanInput and ProductLoanOutput was missing a dash after the
# leading "197" ("197MS-..." -> "197-MS-..."). Fix the text on both
# sheets, then leave the workbook with ProductLoanOutput (B1) selected
# as the active sheet/cell, matching the saved UI state in the commit.

$wb = $excel.ActiveWorkbook

$oldValue = "197MS-EI-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME"
$newValue = "197-MS-EI-DB-DL-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME"

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Correct the product code text on both sheets.
$wsInput.Range("B1").Value  = $newValue
$wsOutput.Range("B1").Value = $newValue

# Restore the input sheet's selection to B1 (it was A6:B6).
$wsInput.Range("B1").Select()

# Make the output sheet the active tab, with B1 selected there too.
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
